# Added Excel Data Map & Hover scenarios
#
# 1) Rename the original sheet "Sheet1" -> "Sheet10"
# 2) Insert a new worksheet "Item Filter" after it, and populate it with
#    the item-filter test data (TC ID / Categories / Size / Color / Availability)
# 3) Re-point the selections on both sheets and make "Item Filter" the active tab

$wb = $excel.ActiveWorkbook

# --- existing sheet: rename + move the selection -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet10"

# --- new sheet -------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Item Filter"

$headerRow = @("TC ID", "Categories", "Size", "Color", "Availability")
for ($c = 1; $c -le $headerRow.Length; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headerRow[$c - 1]
}

$dataRows = @(
    @("tc_01", "Summer Dresses", "L", "Yellow", "In stock"),
    @("tc_02", "Casual Dresses", "S", "Yellow", "In stock"),
    @("tc_03", "Evening Dresses", "M", "Yellow", "In stock")
)

for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $row = $dataRows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Column widths sized to fit the new content (Categories / Availability)
$ws2.Columns.Item(2).ColumnWidth = 15.1666666666667
$ws2.Columns.Item(5).ColumnWidth = 10.3072916666667

# --- selections & active sheet ---------------------------------------------
$ws1.Range("D11").Select()
$ws2.Range("D7").Select()
$ws2.Activate()
